$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 860.3333
$ws.Range("I41").Value = 667.7778
$ws.Range("J41").Value = 1052.8889
$ws.Range("K41").Value = 667.7778
$ws.Range("L41").Value = 1052.8889
$ws.Range("M41").Value = -227.7778
$ws.Range("N41").Value = -1932.8889
$ws.Range("H48").Value = 150
$ws.Range("I48").Value = 150
$ws.Range("K48").Value = 450
$ws.Range("M48").Value = -158
$ws.Range("H53").Value = 1354.1111
$ws.Range("I53").Value = 442.5
$ws.Range("J53").Value = 1614.5714
$ws.Range("K53").Value = 442.5
$ws.Range("L53").Value = 1614.5714
$ws.Range("M53").Value = 194.5
$ws.Range("N53").Value = -2888.5714
$ws.Range("H55").Value = 2374.889
$ws.Range("I55").Value = 185.2
$ws.Range("J55").Value = 5112
$ws.Range("K55").Value = 185.2
$ws.Range("L55").Value = 5112
$ws.Range("M55").Value = 28.80000000000001
$ws.Range("N55").Value = -5540
$ws.Range("H56").Value = 150
$ws.Range("I56").Value = 150
$ws.Range("K56").Value = 450
$ws.Range("M56").Value = 84
$ws.Range("H98").Value = 1387.7826
$ws.Range("I98").Value = 1245.95
$ws.Range("K98").Value = 1245.95
$ws.Range("M98").Value = 252.05
$ws.Range("H107").Value = 19609472
$ws.Range("J107").Value = 958.3333
$ws.Range("L107").Value = 958.3333
$ws.Range("N107").Value = -4798.3333
$ws.Range("H112").Value = 6067.048
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("H122").Value = 1387.7826
$ws.Range("I122").Value = 1245.95
$ws.Range("K122").Value = 3737.85
$ws.Range("M122").Value = -1287.85
$ws.Range("H132").Value = 25002648
$ws.Range("I132").Value = 27780604
$ws.Range("K132").Value = 83341812
$ws.Range("M132").Value = -83339282
$ws.Range("H141").Value = 2302.6365
$ws.Range("I141").Value = 1827
$ws.Range("K141").Value = 5481
$ws.Range("M141").Value = -301
$ws.Range("M112").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 249.16667
$ws.Range("I5").Value = 279
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 279
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = -167
$ws.Range("N5").Value = -324
$ws.Range("H24").Value = 24838.5
$ws.Range("J24").Value = 24838.5
$ws.Range("L24").Value = 24838.5
$ws.Range("N24").Value = -25586.5
$ws.Range("H32").Value = 4163.2856
$ws.Range("I32").Value = 2731.6057
$ws.Range("K32").Value = 2731.6057
$ws.Range("M32").Value = -2444.6057
$ws.Range("H45").Value = 5757206
$ws.Range("I45").Value = 10276947
$ws.Range("K45").Value = 10276947
$ws.Range("M45").Value = -10276570
$ws.Range("H61").Value = 1985.5834
$ws.Range("I61").Value = 1203.6666
$ws.Range("K61").Value = 1203.6666
$ws.Range("M61").Value = -991.6666
$ws.Range("H74").Value = 71891.89999999999
$ws.Range("I74").Value = 5107.9355
$ws.Range("K74").Value = 5107.9355
$ws.Range("M74").Value = -4233.9355
$ws.Range("H77").Value = 71891.89999999999
$ws.Range("I77").Value = 5107.9355
$ws.Range("K77").Value = 25539.6775
$ws.Range("M77").Value = -21171.6775
$ws.Range("H100").Value = 24838.5
$ws.Range("J100").Value = 24838.5
$ws.Range("L100").Value = 24838.5
$ws.Range("N100").Value = -27002.5
$ws.Range("H102").Value = 4389593
$ws.Range("I102").Value = 4905133.5
$ws.Range("J102").Value = 7500
$ws.Range("K102").Value = 4905133.5
$ws.Range("L102").Value = 7500
$ws.Range("M102").Value = -4903511.5
$ws.Range("N102").Value = -10744
$ws.Range("H135").Value = 166708060
$ws.Range("J135").Value = 166708060
$ws.Range("L135").Value = 166708060
$ws.Range("N135").Value = -166718200
$ws.Range("H136").Value = 1985.5834
$ws.Range("I136").Value = 1203.6666
$ws.Range("K136").Value = 3610.9998
$ws.Range("M136").Value = -1060.9998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 249.16667
$ws.Range("I4").Value = 279
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 279
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -164
$ws.Range("N4").Value = -330
$ws.Range("H5").Value = 1150
$ws.Range("I5").Value = 1150
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1150
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1037
$ws.Range("H35").Value = 49998
$ws.Range("J35").Value = 49998
$ws.Range("L35").Value = 49998
$ws.Range("N35").Value = -50618
$ws.Range("N5").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25524.5
$ws.Range("I31").Value = 1105.8462
$ws.Range("K31").Value = 1105.8462
$ws.Range("M31").Value = -810.8462
$ws.Range("H34").Value = 25524.5
$ws.Range("I34").Value = 1105.8462
$ws.Range("K34").Value = 1105.8462
$ws.Range("M34").Value = -903.8462
$ws.Range("H86").Value = 7372.6665
$ws.Range("I86").Value = 6499.647
$ws.Range("J86").Value = 8514.308000000001
$ws.Range("K86").Value = 6499.647
$ws.Range("L86").Value = 8514.308000000001
$ws.Range("M86").Value = -5376.647
$ws.Range("N86").Value = -10760.308
$ws.Range("H89").Value = 7372.6665
$ws.Range("I89").Value = 6499.647
$ws.Range("J89").Value = 8514.308000000001
$ws.Range("K89").Value = 32498.235
$ws.Range("L89").Value = 42571.54000000001
$ws.Range("M89").Value = -26882.235
$ws.Range("N89").Value = -53803.54000000001
$ws.Range("H109").Value = 32597.2
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 32597.2
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 32597.2
$ws.Range("N109").Value = -34677.2
$ws.Range("H134").Value = 2832.6592
$ws.Range("I134").Value = 2242.5588
$ws.Range("K134").Value = 6727.676399999999
$ws.Range("M134").Value = -4192.676399999999
$ws.Range("M109").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1643.3334
$ws.Range("I3").Value = 1643.3334
$ws.Range("K3").Value = 4930.0002
$ws.Range("M3").Value = -4818.0002
$ws.Range("H37").Value = 61159.5
$ws.Range("J37").Value = 61159.5
$ws.Range("L37").Value = 183478.5
$ws.Range("N37").Value = -183702.5
$ws.Range("H98").Value = 2162
$ws.Range("J98").Value = 2500
$ws.Range("L98").Value = 7500
$ws.Range("N98").Value = -10496
$ws.Range("H133").Value = 4247.5
$ws.Range("I133").Value = 4247.5
$ws.Range("K133").Value = 12742.5
$ws.Range("M133").Value = -7682.5
$ws.Range("H136").Value = 2137.5
$ws.Range("I136").Value = 2137.5
$ws.Range("K136").Value = 6412.5
$ws.Range("M136").Value = -1312.5
$ws.Range("H139").Value = 1731.9
$ws.Range("I139").Value = 1731.9
$ws.Range("K139").Value = 5195.700000000001
$ws.Range("M139").Value = -55.70000000000073
$ws.Range("H140").Value = 3959.625
$ws.Range("I140").Value = 3668.2856
$ws.Range("J140").Value = 5999
$ws.Range("K140").Value = 11004.8568
$ws.Range("L140").Value = 17997
$ws.Range("M140").Value = -5824.856800000001
$ws.Range("N140").Value = -28357

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 39995
$ws.Range("J98").Value = 39995
$ws.Range("L98").Value = 39995
$ws.Range("N98").Value = -45985
$ws.Range("H122").Value = 540053.6
$ws.Range("J122").Value = 3698.75
$ws.Range("L122").Value = 11096.25
$ws.Range("N122").Value = -15996.25
$ws.Range("H126").Value = 6442645
$ws.Range("I126").Value = 2844035.5
$ws.Range("J126").Value = 20837082
$ws.Range("K126").Value = 8532106.5
$ws.Range("L126").Value = 62511246
$ws.Range("M126").Value = -8529636.5
$ws.Range("N126").Value = -62516186

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1659.2941
$ws.Range("I16").Value = 1487.2727
$ws.Range("K16").Value = 1487.2727
$ws.Range("M16").Value = -1317.2727
$ws.Range("H22").Value = 222971.75
$ws.Range("I22").Value = 222971.75
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 222971.75
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -222676.75
$ws.Range("H27").Value = 222971.75
$ws.Range("I27").Value = 222971.75
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 222971.75
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -222864.75
$ws.Range("H109").Value = 21988
$ws.Range("J109").Value = 21988
$ws.Range("L109").Value = 21988
$ws.Range("N109").Value = -24762
$ws.Range("H136").Value = 66673.28
$ws.Range("I136").Value = 80290.234
$ws.Range("J136").Value = 7666.5
$ws.Range("K136").Value = 240870.702
$ws.Range("L136").Value = 22999.5
$ws.Range("M136").Value = -238320.702
$ws.Range("N136").Value = -28099.5
$ws.Range("N22").ClearContents()
$ws.Range("N27").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9647.695
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 9647.695
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 9647.695
$ws.Range("N62").Value = -10895.695
$ws.Range("H65").Value = 9647.695
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 9647.695
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 48238.475
$ws.Range("N65").Value = -54478.475
$ws.Range("H107").Value = 35716530
$ws.Range("I107").Value = 45457144
$ws.Range("K107").Value = 136371432
$ws.Range("M107").Value = -136369512
$ws.Range("H132").Value = 25921706
$ws.Range("I132").Value = 27779018
$ws.Range("K132").Value = 83337054
$ws.Range("M132").Value = -83334524
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()
